$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1.52
    3  = 4.13
    4  = 0.15
    5  = 0.1
    6  = 4.55
    7  = 0.48
    8  = 2.3
    9  = 4.2
    10 = 0.027
    11 = 4.25
    12 = 0.11
    13 = 1.52
    14 = 0.25
    15 = 1.05
    16 = 0.1
    17 = 2.1
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Range("I$row").Value = $val
    $ws.Range("Z$row").Value = $val
}
